$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with default (unstyled) formatting, used to strip any
# number-format style stamped onto cells we must force to remain text.
$refStyle = $ws.Range('B2').Style

# --- Plain text / non-numeric-looking values: direct assignment ---
$ws.Range('D2').Value = '42.347.18'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '2.272.54'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  +0.67%  '
$ws.Range('E6').Value = '  +3.07%  '
$ws.Range('E7').Value = '  -0.80%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('E10').Value = '  +1.92%  '
$ws.Range('E11').Value = '  -1.57%  '
$ws.Range('E12').Value = '  +0.11%  '
$ws.Range('E13').Value = '  +3.29%  '
$ws.Range('D14').Value = '2.623.47'
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('E15').Value = '  +2.63%  '
$ws.Range('D16').Value = '2.310.05'
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('E17').Value = '  +0.10%  '
$ws.Range('D18').Value = '42.184.13'
$ws.Range('E18').Value = '  +0.32%  '
$ws.Range('E19').Value = '  -2.64%  '
$ws.Range('D20').Value = '0.0₃0907'
$ws.Range('E20').Value = '  -1.19%  '
$ws.Range('E21').Value = '  +0.65%  '
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('E23').Value = '  -2.09%  '
$ws.Range('E24').Value = '  -0.95%  '
$ws.Range('E25').Value = '  +0.35%  '
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('E27').Value = '  -2.00%  '
$ws.Range('E28').Value = '  +5.21%  '
$ws.Range('E29').Value = '  -1.84%  '
$ws.Range('E30').Value = '  +0.57%  '
$ws.Range('E31').Value = '  +0.37%  '
$ws.Range('E32').Value = '  -1.72%  '
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('E34').Value = '  +3.01%  '
$ws.Range('E35').Value = '  -2.05%  '
$ws.Range('E36').Value = '  +2.89%  '
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('E38').Value = '  -2.78%  '
$ws.Range('E39').Value = '  +0.95%  '
$ws.Range('E40').Value = '  -1.58%  '
$ws.Range('E41').Value = '  -2.58%  '
$ws.Range('E42').Value = '  +2.18%  '
$ws.Range('D43').Value = '1.952.37'
$ws.Range('E43').Value = '  -3.38%  '
$ws.Range('E44').Value = '  -4.30%  '
$ws.Range('E45').Value = '  -0.77%  '
$ws.Range('E46').Value = '  -2.59%  '
$ws.Range('E47').Value = '  -0.82%  '
$ws.Range('E48').Value = '  +0.16%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('E51').Value = '  -1.93%  '

# --- Numeric-looking values that must stay TEXT (match source inlineStr) ---
# Force text entry via NumberFormat "@" so the literal digits (incl. trailing
# zeros) are preserved exactly, then restore the default (unstyled) format.
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '307.20'
$c.Style = $refStyle
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '97.29'
$c.Style = $refStyle
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.527'
$c.Style = $refStyle
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.495'
$c.Style = $refStyle
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '35.35'
$c.Style = $refStyle
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0792'
$c.Style = $refStyle
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '6.90'
$c.Style = $refStyle
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '14.81'
$c.Style = $refStyle
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '0.795'
$c.Style = $refStyle
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '12.48'
$c.Style = $refStyle
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '6.05'
$c.Style = $refStyle
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '68.18'
$c.Style = $refStyle
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '238.74'
$c.Style = $refStyle
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '2.58'
$c.Style = $refStyle
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '1.95'
$c.Style = $refStyle
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '23.66'
$c.Style = $refStyle
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '38.00'
$c.Style = $refStyle
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '9.53'
$c.Style = $refStyle
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '2.12'
$c.Style = $refStyle
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '161.91'
$c.Style = $refStyle
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '5.25'
$c.Style = $refStyle
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = $refStyle
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.0738'
$c.Style = $refStyle
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '17.57'
$c.Style = $refStyle
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '4.09'
$c.Style = $refStyle
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '18.87'
$c.Style = $refStyle
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '9.97'
$c.Style = $refStyle
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '2.90'
$c.Style = $refStyle
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '53.62'
$c.Style = $refStyle
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '92.47'
$c.Style = $refStyle
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '72.05'
$c.Style = $refStyle
